$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.240179
$ws.Range("H2").Value = 3.720537
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.32975033333333
$ws.Range("N2").Value = 48.989251
$ws.Range("O2").Value = 0.3418592381614251
$ws.Range("P2").Value = 0.3418592381614251
$ws.Range("Q2").Value = 20.251813438643
$ws.Range("R2").Value = 182.266320947787
$ws.Range("S2").Value = 0.3418592381614251
$ws.Range("T2").Value = 0.3418592381614251

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.240179
$ws.Range("H3").Value = 3.720537
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.81766766666667
$ws.Range("N3").Value = 80.453003
$ws.Range("O3").Value = 0.5614211638667195
$ws.Range("P3").Value = 0.5614211638667195
$ws.Range("Q3").Value = 33.258708269179
$ws.Range("R3").Value = 299.328374422611
$ws.Range("S3").Value = 0.5614211638667195
$ws.Range("T3").Value = 0.5614211638667195

# Row 4 (Target cluster: sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.240179
$ws.Range("H4").Value = 3.720537
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.620050333333333
$ws.Range("N4").Value = 13.860151
$ws.Range("O4").Value = 0.09671959797185539
$ws.Range("P4").Value = 0.09671959797185539
$ws.Range("Q4").Value = 5.729689402342999
$ws.Range("R4").Value = 51.56720462108699
$ws.Range("S4").Value = 0.09671959797185539
$ws.Range("T4").Value = 0.09671959797185539
